$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.234694480895996
$ws.Range("B1").Value = 3.122582197189331
$ws.Range("C1").Value = 4.148056507110596
$ws.Range("D1").Value = 0.1956574469804764
$ws.Range("E1").Value = 0.236931636929512
